$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "추석연휴 제주 관광객 20만명 넘었다…제주공항 북적"
$ws.Range("B6").Value = "http://www.wowtv.co.kr/NewsCenter/News/Read?articleId=A202010020048&t=NN"
$ws.Range("C6").Value = "한국경제TV"

$ws.Range("A7").Value = "추석인데…모자간에 남매간에 잇단 칼부림 '비극'(종합)"
$ws.Range("B7").Value = "http://news.tf.co.kr/read/national/1816181.htm"
$ws.Range("C7").Value = "더팩트"

$ws.Range("A8").Value = "구름 사이 '휘영청' 추석 보름달, 가장 높게 뜨는 시각은?"
$ws.Range("B8").Value = "https://biz.chosun.com/site/data/html_dir/2020/10/01/2020100100136.html?utm_source=naver&utm_medium=original&utm_campaign=biz"
$ws.Range("C8").Value = "조선비즈"

$ws.Range("A9").Value = "소시지 3개에 김치 3조각... 격리된 병사가 먹은 추석날 저녁"
$ws.Range("B9").Value = "http://www.ohmynews.com/NWS_Web/View/at_pg.aspx?CNTN_CD=A0002680782&CMPT_CD=P0010&utm_source=naver&utm_medium=newsearch&utm_campaign=naver_news"
$ws.Range("C9").Value = "오마이뉴스"

$ws.Range("A10").Value = "이낙연 대표, 추석 맞아 고 노무현 대통령 묘소 참배"
$ws.Range("B10").Value = "https://imnews.imbc.com/news/2020/politics/article/5927450_32626.html"
$ws.Range("C10").Value = "MBC"

$ws.Range("A11").Value = "추미애·조국, 추석날 주거니받거니 “일부 정치검찰 정권과 결탁”"
$ws.Range("B11").Value = "https://www.chosun.com/national/national_general/2020/10/01/7BNGFKXKWZELRL2OONEPG2TKCU/?utm_source=naver&utm_medium=original&utm_campaign=news"
$ws.Range("C11").Value = "조선일보"
